$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sum" column (G), which held =SUM(A:F) formulas, was deleted in the
# edit. Deleting the entire column shifts H->G, I->H, J->I, K->J, L->K,
# which matches every cell-reference / shared-string-index change in the
# recorded diff (shared string "sum" itself is also dropped since it's no
# longer referenced anywhere, which Excel does automatically on save).
$ws.Range("G1:G1048576").Select()
$ws.Range("G1:G1048576").EntireColumn.Delete()

# Reproduce the final recorded cursor position (row 15 of what is now
# column G, formerly column H).
$ws.Range("G15").Select()
